# Update Lgi2-Adam22 LR-pairs sheet with recomputed TPM-based NATMI statistics.
# The workbook has a single worksheet ("Sheet1"); update the affected data cells
# (columns E-T, rows 2-19) in place with the newly computed values from the pipeline re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 5).Value = 3  # E2: 1 -> 3
$ws.Cells.Item(2, 6).Value = 1  # F2: 0.3333333333333333 -> 1
$ws.Cells.Item(2, 7).Value = 0.1376636666666667  # G2: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(2, 8).Value = 0.412991  # H2: 0.010366 -> 0.412991
$ws.Cells.Item(2, 9).Value = 0.01821680097623009  # I2: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(2, 10).Value = 0.01821680097623009  # J2: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(2, 13).Value = 1.667434  # M2: 1.715363 -> 1.667434
$ws.Cells.Item(2, 14).Value = 5.002302  # N2: 5.146089 -> 5.002302
$ws.Cells.Item(2, 15).Value = 0.3223739883484499  # O2: 0.1765901754351529 -> 0.3223739883484499
$ws.Cells.Item(2, 16).Value = 0.32237398834845  # P2: 0.1765901754351529 -> 0.32237398834845
$ws.Cells.Item(2, 17).Value = 0.2295450783646667  # Q2: 0.005927150952666667 -> 0.2295450783646667
$ws.Cells.Item(2, 18).Value = 2.065905705282  # R2: 0.053344358574 -> 2.065905705282
$ws.Cells.Item(2, 19).Value = 0.005872622785657229  # S2: 7.93155649478202E-05 -> 0.005872622785657229
$ws.Cells.Item(2, 20).Value = 0.005872622785657231  # T2: 7.931556494782017E-05 -> 0.005872622785657231

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 5).Value = 3  # E3: 1 -> 3
$ws.Cells.Item(3, 6).Value = 1  # F3: 0.3333333333333333 -> 1
$ws.Cells.Item(3, 7).Value = 0.1376636666666667  # G3: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(3, 8).Value = 0.412991  # H3: 0.010366 -> 0.412991
$ws.Cells.Item(3, 9).Value = 0.01821680097623009  # I3: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(3, 10).Value = 0.01821680097623009  # J3: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(3, 15).Value = 0.2193354457157105  # O3: 0.1167906248092277 -> 0.2193354457157105
$ws.Cells.Item(3, 16).Value = 0.2193354457157106  # P3: 0.1167906248092277 -> 0.2193354457157106
$ws.Cells.Item(3, 17).Value = 0.1561769059994444  # Q3: 0.003920012318888889 -> 0.1561769059994444
$ws.Cells.Item(3, 18).Value = 1.405592153995  # R3: 0.03528011087 -> 1.405592153995
$ws.Cells.Item(3, 19).Value = 0.003995590161635816  # S3: 5.245656710247993E-05 -> 0.003995590161635816
$ws.Cells.Item(3, 20).Value = 0.003995590161635818  # T3: 5.245656710247993E-05 -> 0.003995590161635818

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 5).Value = 3  # E4: 1 -> 3
$ws.Cells.Item(4, 6).Value = 1  # F4: 0.3333333333333333 -> 1
$ws.Cells.Item(4, 7).Value = 0.1376636666666667  # G4: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(4, 8).Value = 0.412991  # H4: 0.010366 -> 0.412991
$ws.Cells.Item(4, 9).Value = 0.01821680097623009  # I4: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(4, 10).Value = 0.01821680097623009  # J4: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(4, 13).Value = 0.2055123333333333  # M4: 0.8505606666666666 -> 0.2055123333333333
$ws.Cells.Item(4, 14).Value = 0.616537  # N4: 2.551682 -> 0.616537
$ws.Cells.Item(4, 15).Value = 0.03973280534729575  # O4: 0.08756202468218523 -> 0.03973280534729575
$ws.Cells.Item(4, 16).Value = 0.03973280534729576  # P4: 0.08756202468218521 -> 0.03973280534729576
$ws.Cells.Item(4, 17).Value = 0.02829158135188889  # Q4: 0.002938970623555555 -> 0.02829158135188889
$ws.Cells.Item(4, 18).Value = 0.254624232167  # R4: 0.026450735612 -> 0.254624232167
$ws.Cells.Item(4, 19).Value = 0.0007238046072389773  # S4: 3.932852684770584E-05 -> 0.0007238046072389773
$ws.Cells.Item(4, 20).Value = 0.0007238046072389775  # T4: 3.932852684770583E-05 -> 0.0007238046072389775

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 5).Value = 3  # E5: 1 -> 3
$ws.Cells.Item(5, 6).Value = 1  # F5: 0.3333333333333333 -> 1
$ws.Cells.Item(5, 7).Value = 0.1376636666666667  # G5: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(5, 8).Value = 0.412991  # H5: 0.010366 -> 0.412991
$ws.Cells.Item(5, 9).Value = 0.01821680097623009  # I5: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(5, 10).Value = 0.01821680097623009  # J5: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(5, 13).Value = 1.661741333333333  # M5: 3.547937 -> 1.661741333333333
$ws.Cells.Item(5, 14).Value = 4.985224  # N5: 10.643811 -> 4.985224
$ws.Cells.Item(5, 15).Value = 0.3212733944672698  # O5: 0.3652467829041842 -> 0.3212733944672698
$ws.Cells.Item(5, 16).Value = 0.3212733944672699  # P5: 0.3652467829041842 -> 0.3212733944672699
$ws.Cells.Item(5, 17).Value = 0.2287614049982222  # Q5: 0.01225930498066666 -> 0.2287614049982222
$ws.Cells.Item(5, 18).Value = 2.058852644984  # R5: 0.110333744826 -> 2.058852644984
$ws.Cells.Item(5, 19).Value = 0.005852573485968115  # S5: 0.0001640507738328705 -> 0.005852573485968115
$ws.Cells.Item(5, 20).Value = 0.005852573485968116  # T5: 0.0001640507738328705 -> 0.005852573485968116

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6, 5).Value = 3  # E6: 1 -> 3
$ws.Cells.Item(6, 6).Value = 1  # F6: 0.3333333333333333 -> 1
$ws.Cells.Item(6, 7).Value = 0.1376636666666667  # G6: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(6, 8).Value = 0.412991  # H6: 0.010366 -> 0.412991
$ws.Cells.Item(6, 9).Value = 0.01821680097623009  # I6: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(6, 10).Value = 0.01821680097623009  # J6: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(6, 13).Value = 0.1178836666666667  # M6: 0.1715316666666667 -> 0.1178836666666667
$ws.Cells.Item(6, 14).Value = 0.353651  # N6: 0.514595 -> 0.353651
$ws.Cells.Item(6, 15).Value = 0.02279108365576842  # O6: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(6, 16).Value = 0.02279108365576842  # P6: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(6, 17).Value = 0.01622829779344444  # Q6: 0.0005926990855555555 -> 0.01622829779344444
$ws.Cells.Item(6, 18).Value = 0.146054680141  # R6: 0.00533429177 -> 0.146054680141
$ws.Cells.Item(6, 19).Value = 0.0004151806349897437  # S6: 7.931342257066197E-06 -> 0.0004151806349897437
$ws.Cells.Item(6, 20).Value = 0.0004151806349897437  # T6: 7.931342257066197E-06 -> 0.0004151806349897437

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7, 5).Value = 3  # E7: 1 -> 3
$ws.Cells.Item(7, 6).Value = 1  # F7: 0.3333333333333333 -> 1
$ws.Cells.Item(7, 7).Value = 0.1376636666666667  # G7: 0.003455333333333333 -> 0.1376636666666667
$ws.Cells.Item(7, 8).Value = 0.412991  # H7: 0.010366 -> 0.412991
$ws.Cells.Item(7, 9).Value = 0.01821680097623009  # I7: 0.0004491504963533279 -> 0.01821680097623009
$ws.Cells.Item(7, 10).Value = 0.01821680097623009  # J7: 0.0004491504963533278 -> 0.01821680097623009
$ws.Cells.Item(7, 13).Value = 0.385306  # M7: 2.293933666666666 -> 0.385306
$ws.Cells.Item(7, 14).Value = 1.155918  # N7: 6.881800999999999 -> 1.155918
$ws.Cells.Item(7, 15).Value = 0.07449328246550557  # O7: 0.2361518516099917 -> 0.07449328246550557
$ws.Cells.Item(7, 16).Value = 0.0744932824655056  # P7: 0.2361518516099917 -> 0.0744932824655056
$ws.Cells.Item(7, 17).Value = 0.05304263674866666  # Q7: 0.007926305462888887 -> 0.05304263674866666
$ws.Cells.Item(7, 18).Value = 0.477383730738  # R7: 0.071336749166 -> 0.477383730738
$ws.Cells.Item(7, 19).Value = 0.001357029300740206  # S7: 0.0001060677213653852 -> 0.001357029300740206
$ws.Cells.Item(7, 20).Value = 0.001357029300740206  # T7: 0.0001060677213653852 -> 0.001357029300740206

# Row 8: FAPs -> ECs
$ws.Cells.Item(8, 9).Value = 0.9679346439276632  # I8: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(8, 10).Value = 0.967934643927663  # J8: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(8, 13).Value = 1.667434  # M8: 1.715363 -> 1.667434
$ws.Cells.Item(8, 14).Value = 5.002302  # N8: 5.146089 -> 5.002302
$ws.Cells.Item(8, 15).Value = 0.3223739883484499  # O8: 0.1765901754351529 -> 0.3223739883484499
$ws.Cells.Item(8, 16).Value = 0.32237398834845  # P8: 0.1765901754351529 -> 0.32237398834845
$ws.Cells.Item(8, 17).Value = 12.19668777093  # Q8: 12.547271391135 -> 12.19668777093
$ws.Cells.Item(8, 18).Value = 109.77018993837  # R8: 112.925442520215 -> 109.77018993837
$ws.Cells.Item(8, 19).Value = 0.3120369516235975  # S8: 0.1679042641041139 -> 0.3120369516235975
$ws.Cells.Item(8, 20).Value = 0.3120369516235975  # T8: 0.1679042641041138 -> 0.3120369516235975

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9, 9).Value = 0.9679346439276632  # I9: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(9, 10).Value = 0.967934643927663  # J9: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(9, 15).Value = 0.2193354457157105  # O9: 0.1167906248092277 -> 0.2193354457157105
$ws.Cells.Item(9, 16).Value = 0.2193354457157106  # P9: 0.1167906248092277 -> 0.2193354457157106
$ws.Cells.Item(9, 19).Value = 0.2123023765495515  # S9: 0.1110460639417285 -> 0.2123023765495515
$ws.Cells.Item(9, 20).Value = 0.2123023765495516  # T9: 0.1110460639417285 -> 0.2123023765495516

# Row 10: FAPs -> Inflammatory-Mac
$ws.Cells.Item(10, 9).Value = 0.9679346439276632  # I10: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(10, 10).Value = 0.967934643927663  # J10: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(10, 13).Value = 0.2055123333333333  # M10: 0.8505606666666666 -> 0.2055123333333333
$ws.Cells.Item(10, 14).Value = 0.616537  # N10: 2.551682 -> 0.616537
$ws.Cells.Item(10, 15).Value = 0.03973280534729575  # O10: 0.08756202468218523 -> 0.03973280534729575
$ws.Cells.Item(10, 16).Value = 0.03973280534729576  # P10: 0.08756202468218521 -> 0.03973280534729576
$ws.Cells.Item(10, 17).Value = 1.503249761455  # Q10: 6.221549327629998 -> 1.503249761455
$ws.Cells.Item(10, 18).Value = 13.529247853095  # R10: 55.99394394866999 -> 13.529247853095
$ws.Cells.Item(10, 19).Value = 0.03845875879608186  # S10: 0.08325512606519504 -> 0.03845875879608186
$ws.Cells.Item(10, 20).Value = 0.03845875879608187  # T10: 0.08325512606519501 -> 0.03845875879608187

# Row 11: FAPs -> MuSCs
$ws.Cells.Item(11, 9).Value = 0.9679346439276632  # I11: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(11, 10).Value = 0.967934643927663  # J11: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(11, 13).Value = 1.661741333333333  # M11: 3.547937 -> 1.661741333333333
$ws.Cells.Item(11, 14).Value = 4.985224  # N11: 10.643811 -> 4.985224
$ws.Cells.Item(11, 15).Value = 0.3212733944672698  # O11: 0.3652467829041842 -> 0.3212733944672698
$ws.Cells.Item(11, 16).Value = 0.3212733944672699  # P11: 0.3652467829041842 -> 0.3212733944672699
$ws.Cells.Item(11, 17).Value = 12.15504793516  # Q11: 25.95189963736499 -> 12.15504793516
$ws.Cells.Item(11, 18).Value = 109.39543141644  # R11: 233.567096736285 -> 109.39543141644
$ws.Cells.Item(11, 19).Value = 0.3109716486771085  # S11: 0.3472814506741473 -> 0.3109716486771085
$ws.Cells.Item(11, 20).Value = 0.3109716486771085  # T11: 0.3472814506741472 -> 0.3109716486771085

# Row 12: FAPs -> Neutrophils
$ws.Cells.Item(12, 9).Value = 0.9679346439276632  # I12: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(12, 10).Value = 0.967934643927663  # J12: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(12, 13).Value = 0.1178836666666667  # M12: 0.1715316666666667 -> 0.1178836666666667
$ws.Cells.Item(12, 14).Value = 0.353651  # N12: 0.514595 -> 0.353651
$ws.Cells.Item(12, 15).Value = 0.02279108365576842  # O12: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(12, 16).Value = 0.02279108365576842  # P12: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(12, 17).Value = 0.8622771729649998  # Q12: 1.254693247925 -> 0.8622771729649998
$ws.Cells.Item(12, 18).Value = 7.760494556684998  # R12: 11.292239231325 -> 7.760494556684998
$ws.Cells.Item(12, 19).Value = 0.02206027944307179  # S12: 0.01678997288749893 -> 0.02206027944307179
$ws.Cells.Item(12, 20).Value = 0.02206027944307179  # T12: 0.01678997288749893 -> 0.02206027944307179

# Row 13: FAPs -> Resolving-Mac
$ws.Cells.Item(13, 9).Value = 0.9679346439276632  # I13: 0.9508131677788118 -> 0.9679346439276632
$ws.Cells.Item(13, 10).Value = 0.967934643927663  # J13: 0.9508131677788116 -> 0.967934643927663
$ws.Cells.Item(13, 13).Value = 0.385306  # M13: 2.293933666666666 -> 0.385306
$ws.Cells.Item(13, 14).Value = 1.155918  # N13: 6.881800999999999 -> 1.155918
$ws.Cells.Item(13, 15).Value = 0.07449328246550557  # O13: 0.2361518516099917 -> 0.07449328246550557
$ws.Cells.Item(13, 16).Value = 0.0744932824655056  # P13: 0.2361518516099917 -> 0.0744932824655056
$ws.Cells.Item(13, 17).Value = 2.818376606369999  # Q13: 16.77931042521499 -> 2.818376606369999
$ws.Cells.Item(13, 18).Value = 25.36538945733  # R13: 151.013793826935 -> 25.36538945733
$ws.Cells.Item(13, 19).Value = 0.07210462883825197  # S13: 0.2245362901061281 -> 0.07210462883825197
$ws.Cells.Item(13, 20).Value = 0.07210462883825199  # T13: 0.2245362901061281 -> 0.07210462883825199

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14, 7).Value = 0.104653  # G14: 0.374941 -> 0.104653
$ws.Cells.Item(14, 8).Value = 0.313959  # H14: 1.124823 -> 0.313959
$ws.Cells.Item(14, 9).Value = 0.01384855509610675  # I14: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(14, 10).Value = 0.01384855509610675  # J14: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(14, 13).Value = 1.667434  # M14: 1.715363 -> 1.667434
$ws.Cells.Item(14, 14).Value = 5.002302  # N14: 5.146089 -> 5.002302
$ws.Cells.Item(14, 15).Value = 0.3223739883484499  # O14: 0.1765901754351529 -> 0.3223739883484499
$ws.Cells.Item(14, 16).Value = 0.32237398834845  # P14: 0.1765901754351529 -> 0.32237398834845
$ws.Cells.Item(14, 17).Value = 0.174501970402  # Q14: 0.6431599185829999 -> 0.174501970402
$ws.Cells.Item(14, 18).Value = 1.570517733618  # R14: 5.788439267246999 -> 1.570517733618
$ws.Cells.Item(14, 19).Value = 0.004464413939195184  # S14: 0.008606595766091255 -> 0.004464413939195184
$ws.Cells.Item(14, 20).Value = 0.004464413939195185  # T14: 0.008606595766091253 -> 0.004464413939195185

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15, 7).Value = 0.104653  # G15: 0.374941 -> 0.104653
$ws.Cells.Item(15, 8).Value = 0.313959  # H15: 1.124823 -> 0.313959
$ws.Cells.Item(15, 9).Value = 0.01384855509610675  # I15: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(15, 10).Value = 0.01384855509610675  # J15: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(15, 15).Value = 0.2193354457157105  # O15: 0.1167906248092277 -> 0.2193354457157105
$ws.Cells.Item(15, 16).Value = 0.2193354457157106  # P15: 0.1167906248092277 -> 0.2193354457157106
$ws.Cells.Item(15, 17).Value = 0.1187269098616667  # Q15: 0.4253636905816666 -> 0.1187269098616667
$ws.Cells.Item(15, 18).Value = 1.068542188755  # R15: 3.828273215235 -> 1.068542188755
$ws.Cells.Item(15, 19).Value = 0.003037479004523148  # S15: 0.005692104300396757 -> 0.003037479004523148
$ws.Cells.Item(15, 20).Value = 0.003037479004523148  # T15: 0.005692104300396757 -> 0.003037479004523148

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16, 7).Value = 0.104653  # G16: 0.374941 -> 0.104653
$ws.Cells.Item(16, 8).Value = 0.313959  # H16: 1.124823 -> 0.313959
$ws.Cells.Item(16, 9).Value = 0.01384855509610675  # I16: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(16, 10).Value = 0.01384855509610675  # J16: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(16, 13).Value = 0.2055123333333333  # M16: 0.8505606666666666 -> 0.2055123333333333
$ws.Cells.Item(16, 14).Value = 0.616537  # N16: 2.551682 -> 0.616537
$ws.Cells.Item(16, 15).Value = 0.03973280534729575  # O16: 0.08756202468218523 -> 0.03973280534729575
$ws.Cells.Item(16, 16).Value = 0.03973280534729576  # P16: 0.08756202468218521 -> 0.03973280534729576
$ws.Cells.Item(16, 17).Value = 0.02150748222033333  # Q16: 0.3189100669206666 -> 0.02150748222033333
$ws.Cells.Item(16, 18).Value = 0.193567339983  # R16: 2.870190602286 -> 0.193567339983
$ws.Cells.Item(16, 19).Value = 0.0005502419439749101  # S16: 0.004267570090142487 -> 0.0005502419439749101
$ws.Cells.Item(16, 20).Value = 0.0005502419439749102  # T16: 0.004267570090142486 -> 0.0005502419439749102

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17, 7).Value = 0.104653  # G17: 0.374941 -> 0.104653
$ws.Cells.Item(17, 8).Value = 0.313959  # H17: 1.124823 -> 0.313959
$ws.Cells.Item(17, 9).Value = 0.01384855509610675  # I17: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(17, 10).Value = 0.01384855509610675  # J17: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(17, 13).Value = 1.661741333333333  # M17: 3.547937 -> 1.661741333333333
$ws.Cells.Item(17, 14).Value = 4.985224  # N17: 10.643811 -> 4.985224
$ws.Cells.Item(17, 15).Value = 0.3212733944672698  # O17: 0.3652467829041842 -> 0.3212733944672698
$ws.Cells.Item(17, 16).Value = 0.3212733944672699  # P17: 0.3652467829041842 -> 0.3212733944672699
$ws.Cells.Item(17, 17).Value = 0.1739062157573333  # Q17: 1.330267046717 -> 0.1739062157573333
$ws.Cells.Item(17, 18).Value = 1.565155941816  # R17: 11.972403420453 -> 1.565155941816
$ws.Cells.Item(17, 19).Value = 0.004449172304193224  # S17: 0.01780128145620402 -> 0.004449172304193224
$ws.Cells.Item(17, 20).Value = 0.004449172304193224  # T17: 0.01780128145620402 -> 0.004449172304193224

# Row 18: MuSCs -> Neutrophils
$ws.Cells.Item(18, 7).Value = 0.104653  # G18: 0.374941 -> 0.104653
$ws.Cells.Item(18, 8).Value = 0.313959  # H18: 1.124823 -> 0.313959
$ws.Cells.Item(18, 9).Value = 0.01384855509610675  # I18: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(18, 10).Value = 0.01384855509610675  # J18: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(18, 13).Value = 0.1178836666666667  # M18: 0.1715316666666667 -> 0.1178836666666667
$ws.Cells.Item(18, 14).Value = 0.353651  # N18: 0.514595 -> 0.353651
$ws.Cells.Item(18, 15).Value = 0.02279108365576842  # O18: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(18, 16).Value = 0.02279108365576842  # P18: 0.01765854055925821 -> 0.02279108365576842
$ws.Cells.Item(18, 17).Value = 0.01233687936766667  # Q18: 0.06431425463166666 -> 0.01233687936766667
$ws.Cells.Item(18, 18).Value = 0.111031914309  # R18: 0.5788282916849999 -> 0.111031914309
$ws.Cells.Item(18, 19).Value = 0.000315623577706887  # S18: 0.0008606363295022159 -> 0.000315623577706887
$ws.Cells.Item(18, 20).Value = 0.000315623577706887  # T18: 0.0008606363295022159 -> 0.000315623577706887

# Row 19: MuSCs -> Resolving-Mac
$ws.Cells.Item(19, 7).Value = 0.104653  # G19: 0.374941 -> 0.104653
$ws.Cells.Item(19, 8).Value = 0.313959  # H19: 1.124823 -> 0.313959
$ws.Cells.Item(19, 9).Value = 0.01384855509610675  # I19: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(19, 10).Value = 0.01384855509610675  # J19: 0.04873768172483497 -> 0.01384855509610675
$ws.Cells.Item(19, 13).Value = 0.385306  # M19: 2.293933666666666 -> 0.385306
$ws.Cells.Item(19, 14).Value = 1.155918  # N19: 6.881800999999999 -> 1.155918
$ws.Cells.Item(19, 15).Value = 0.07449328246550557  # O19: 0.2361518516099917 -> 0.07449328246550557
$ws.Cells.Item(19, 16).Value = 0.0744932824655056  # P19: 0.2361518516099917 -> 0.0744932824655056
$ws.Cells.Item(19, 17).Value = 0.040323428818  # Q19: 0.8600897829136664 -> 0.040323428818
$ws.Cells.Item(19, 18).Value = 0.362910859362  # R19: 7.740808046222998 -> 0.362910859362
$ws.Cells.Item(19, 19).Value = 0.001031624326513397  # S19: 0.01150949378249823 -> 0.001031624326513397
$ws.Cells.Item(19, 20).Value = 0.001031624326513397  # T19: 0.01150949378249823 -> 0.001031624326513397
